# Updated symbol list on Wed Dec 21 09:08:00 UTC 2022 with GitHub Actions
#
# cryptos.xlsx, sheet1: refresh the "Price" (column D) snapshot for the
# rows whose quote moved, and bump every row's "Hora" (column G) from the
# previous scrape hour (8) to the new one (9).
#
# Column D and column G are stored as literal text in the workbook (prices
# like "0.05700" / "--" need their exact text, not a coerced Double), so
# each touched cell is first forced to Text format before the new literal
# is written — otherwise Excel's usual Value auto-typing would turn
# "0.05700" into the number 0.057 and drop the significant trailing zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new Price (column D) text, only for rows whose price changed
$prices = @{
    2  = "248.74"
    3  = "22.45"
    4  = "5.396"
    5  = "0.05700"
    6  = "3.403"
    7  = "6.323"
    8  = "0.8086"
    9  = "0.9232"
    10 = "0.1411"
    11 = "0.07418"
    12 = "0.03113"
    13 = "0.03017"
    14 = "0.09382"
    15 = "3.872"
    16 = "0.001579"
    17 = "0.04744"
    18 = "0.01825"
    19 = "0.0005847"
    20 = "0.006481"
    21 = "0.004993"
    22 = "0.001006"
    23 = "0.0001498"
    24 = "3.699"
    25 = "2.200"
    26 = "0.3257"
    40 = "0.03999"
    41 = "0.006851"
    42 = "0.1070"
    43 = "0.002706"
    44 = "0.007511"
    45 = "0.00005798"
    48 = "0.2088"
    49 = "0.00002099"
}

foreach ($row in $prices.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $prices[$row]
}

# every data row (2-51): Hora (column G) goes from "8" to "9"
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Range("G$row")
    $cell.NumberFormat = "@"
    $cell.Value = "9"
}
